# Insert a new weekly price record for Espinaca (Vega Modelo de Temuco) as
# row 189, pushing the existing rows 189-314 down to 190-315.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(189).Insert()

$ws.Range("A189").Value = 10
$ws.Range("B189").Value = "Vega Modelo de Temuco"
$ws.Range("C189").Value = "La Araucanía"
$ws.Range("D189").Value = 45216
$ws.Range("E189").Value = 9
$ws.Range("F189").Value = 100112012
$ws.Range("G189").Value = "Espinaca"
$ws.Range("H189").Value = "Sin especificar"
$ws.Range("I189").Value = "Primera"
$ws.Range("J189").Value = 125
$ws.Range("K189").Value = 10000
$ws.Range("L189").Value = 10000
$ws.Range("M189").Value = 10000
$ws.Range("N189").Value = "`$/docena de atados"
$ws.Range("O189").Value = "Región de La Araucanía"
$ws.Range("P189").Value = 3333
$ws.Range("Q189").Value = 3
$ws.Range("R189").Value = "Hortaliza"
